$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets, in order, after the existing "Login" sheet.
# ---------------------------------------------------------------------------
$login = $wb.Worksheets.Item(1)
$newAppt = $wb.Worksheets.Add($null, $login)
$newAppt.Name = "New appointment"
$rejectApp = $wb.Worksheets.Add($null, $newAppt)
$rejectApp.Name = "Reject app by Interpreter"

# ---------------------------------------------------------------------------
# 2. "New appointment" sheet content.
#    The fill order below reproduces the order in which the shared strings
#    were first introduced in the authored workbook.
# ---------------------------------------------------------------------------

# Row 1 headers (Department/J1 is deliberately filled in later, to match
# the original authoring order).
$newAppt.Range("A1").Value = "Email Address"
$newAppt.Range("C1").Value = "Appointment Date"
$newAppt.Range("D1").Value = "App Start time"
$newAppt.Range("E1").Value = "App End Time"
$newAppt.Range("F1").Value = "Client"
$newAppt.Range("G1").Value = "Facility"
$newAppt.Range("H1").Value = "App Type"
$newAppt.Range("I1").Value = "Building"
$newAppt.Range("K1").Value = "Patient Name"
$newAppt.Range("L1").Value = "Requested Language"
$newAppt.Range("B1").Value = "Password"

# Row 2 sample data (simple text values first).
$newAppt.Range("F2").Value = "CHOP"
$newAppt.Range("G2").Value = "CHOP Main"
$newAppt.Range("H2").Value = "Fitness"
$newAppt.Range("I2").Value = "Abramson Building"
$newAppt.Range("J2").Value = "CS"
$newAppt.Range("K2").Value = "Harsha"
$newAppt.Range("L2").Value = "spanish"

# Department header filled in afterwards.
$newAppt.Range("J1").Value = "Department"

# Login credentials (reuse existing shared strings).
$newAppt.Range("A2").Value = "ravi.thota@sstech.us"
$newAppt.Range("B2").Value = "Welcome@1"

# Date / time text values added last.
$newAppt.Range("D2").Value = "`"16:15:00`""
$newAppt.Range("E2").Value = "`"16:35:00`""
$newAppt.Range("C2").Value = "`"23-12-2022`""

# Number formats for the date/time text cells.
$newAppt.Range("C2").NumberFormat = "mm-dd-yy"
$newAppt.Range("D2").NumberFormat = "h:mm"
$newAppt.Range("E2").NumberFormat = "h:mm"

# Header row styling: bold-free default font, light accent fill and thin box
# border around every header cell.
$headerRange = $newAppt.Range("A1:L1")
$headerRange.Interior.ThemeColor = 6
$headerRange.Interior.TintAndShade = 0.79998168889431442
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Hyperlinks for the credential cells.
[void]$newAppt.Hyperlinks.Add($newAppt.Range("A2"), "mailto:ravi.thota@sstech.us")
[void]$newAppt.Hyperlinks.Add($newAppt.Range("B2"), "mailto:Welcome@1")

# Column widths.
$newAppt.Columns.Item(1).ColumnWidth = 18.666666666666668
$newAppt.Columns.Item(2).ColumnWidth = 12.0
$newAppt.Columns.Item(3).ColumnWidth = 15.666666666666668
$newAppt.Columns.Item(4).ColumnWidth = 14.666666666666668
$newAppt.Columns.Item(5).ColumnWidth = 16.666666666666668
$newAppt.Columns.Item(7).ColumnWidth = 10.333333333333334
$newAppt.Columns.Item(8).ColumnWidth = 10.5
$newAppt.Columns.Item(9).ColumnWidth = 16.666666666666668
$newAppt.Columns.Item(10).ColumnWidth = 11.0
$newAppt.Columns.Item(11).ColumnWidth = 13.333333333333334
$newAppt.Columns.Item(12).ColumnWidth = 19.0

# ---------------------------------------------------------------------------
# 3. "Reject app by Interpreter" sheet content.
# ---------------------------------------------------------------------------

# Row 1 headers (order follows original authoring: A, B, D, E, then C).
$rejectApp.Range("A1").Value = "Scheduler Username"
$rejectApp.Range("B1").Value = "Scheduler Password"
$rejectApp.Range("D1").Value = "Interpreter Username"
$rejectApp.Range("E1").Value = "Interpreter Password"
$rejectApp.Range("C1").Value = "Interpreter Name"
$rejectApp.Range("F1").Value = "Requested Language"

# Row 2 values.
$rejectApp.Range("C2").Value = "Matt Laborde"
$rejectApp.Range("D2").Value = "matt.laborde@sstech.us"
$rejectApp.Range("F2").Value = "Spanish"
$rejectApp.Range("A2").Value = "ravi.thota@sstech.us"
$rejectApp.Range("B2").Value = "Welcome@1"
$rejectApp.Range("E2").Value = "Welcome@1"

# Hyperlinks (order matters for the relationship id assignment).
[void]$rejectApp.Hyperlinks.Add($rejectApp.Range("A2"), "mailto:ravi.thota@sstech.us")
[void]$rejectApp.Hyperlinks.Add($rejectApp.Range("B2"), "mailto:Welcome@1")
[void]$rejectApp.Hyperlinks.Add($rejectApp.Range("E2"), "mailto:Welcome@1")
[void]$rejectApp.Hyperlinks.Add($rejectApp.Range("D2"), "mailto:matt.laborde@sstech.us")

# Column widths.
$rejectApp.Columns.Item(1).ColumnWidth = 17.333333333333336
$rejectApp.Columns.Item(2).ColumnWidth = 18.666666666666664
$rejectApp.Columns.Item(3).ColumnWidth = 18.666666666666664
$rejectApp.Columns.Item(4).ColumnWidth = 30.166666666666664
$rejectApp.Columns.Item(5).ColumnWidth = 18.333333333333336
$rejectApp.Columns.Item(6).ColumnWidth = 23.666666666666664

# ---------------------------------------------------------------------------
# 4. Selections / active sheet bookkeeping, matching the authored workbook:
#    "New appointment" ends up the active (selected) tab, with C2 selected;
#    "Reject app by Interpreter" keeps a pending selection of F7.
# ---------------------------------------------------------------------------
$rejectApp.Activate()
[void]$rejectApp.Range("F7").Select()

$newAppt.Activate()
[void]$newAppt.Range("C2").Select()
